# Rename the worksheet from "nucleotide_table" to "modification_sets"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "modification_sets"

# Move the active selection to cell A36
$ws.Range("A36").Select()
